# chore: update Sheets via scheduled runner
# Applies numeric value updates (and a few cell adds/removes) to the
# per-profession Leve profit tables (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# to reflect refreshed market-board price data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 685.68964
$ws.Range("I17").Value = 280
$ws.Range("J17").Value = 707.8182
$ws.Range("K17").Value = 840
$ws.Range("L17").Value = 2123.4546
$ws.Range("M17").Value = -672
$ws.Range("N17").Value = -2459.4546

$ws.Range("H18").Value = 297
$ws.Range("I18").Value = 297
$ws.Range("K18").Value = 297
$ws.Range("M18").Value = -13

$ws.Range("H26").Value = 49500
$ws.Range("J26").Value = 49000
$ws.Range("L26").Value = 49000
$ws.Range("N26").Value = -49688

$ws.Range("H76").Value = 5053587
$ws.Range("I76").Value = 5294053
$ws.Range("J76").Value = 3800
$ws.Range("K76").Value = 5294053
$ws.Range("L76").Value = 3800
$ws.Range("M76").Value = -5293738
$ws.Range("N76").Value = -4430

$ws.Range("H79").Value = 5053587
$ws.Range("I79").Value = 5294053
$ws.Range("J79").Value = 3800
$ws.Range("K79").Value = 5294053
$ws.Range("L79").Value = 3800
$ws.Range("M79").Value = -5292961
$ws.Range("N79").Value = -5984

$ws.Range("H118").Value = 900
$ws.Range("I118").Value = 800
$ws.Range("K118").Value = 2400
$ws.Range("M118").Value = -743

$ws.Range("H129").Value = 1088.3667
$ws.Range("J129").Value = 1166.4445
$ws.Range("L129").Value = 3499.3335
$ws.Range("N129").Value = -13499.3335

$ws.Range("H137").Value = 1714.762
$ws.Range("I137").Value = 1050
$ws.Range("J137").Value = 2047.1428
$ws.Range("K137").Value = 3150
$ws.Range("L137").Value = 6141.428400000001
$ws.Range("M137").Value = -600
$ws.Range("N137").Value = -11241.4284

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4025.111
$ws.Range("I32").Value = 3065.6897
$ws.Range("J32").Value = 7999.857
$ws.Range("K32").Value = 3065.6897
$ws.Range("L32").Value = 7999.857
$ws.Range("M32").Value = -2778.6897
$ws.Range("N32").Value = -8573.857

$ws.Range("H63").Value = 3639.5454
$ws.Range("I63").Value = 3558.6843
$ws.Range("J63").Value = 4151.6665
$ws.Range("K63").Value = 3558.6843
$ws.Range("L63").Value = 4151.6665
$ws.Range("M63").Value = -2872.6843
$ws.Range("N63").Value = -5523.6665

$ws.Range("H66").Value = 3639.5454
$ws.Range("I66").Value = 3558.6843
$ws.Range("J66").Value = 4151.6665
$ws.Range("K66").Value = 17793.4215
$ws.Range("L66").Value = 20758.3325
$ws.Range("M66").Value = -14361.4215
$ws.Range("N66").Value = -27622.3325

$ws.Range("H88").Value = 4912.5
$ws.Range("I88").Value = 3000
$ws.Range("J88").Value = 6060
$ws.Range("K88").Value = 3000
$ws.Range("L88").Value = 6060
$ws.Range("M88").Value = -2594
$ws.Range("N88").Value = -6872

$ws.Range("H91").Value = 4912.5
$ws.Range("I91").Value = 3000
$ws.Range("J91").Value = 6060
$ws.Range("K91").Value = 3000
$ws.Range("L91").Value = 6060
$ws.Range("M91").Value = -1596
$ws.Range("N91").Value = -8868

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 14476.75
$ws.Range("I86").Value = 2133.3333
$ws.Range("J86").Value = 21882.8
$ws.Range("K86").Value = 2133.3333
$ws.Range("L86").Value = 21882.8
$ws.Range("M86").Value = -1010.3333
$ws.Range("N86").Value = -24128.8

$ws.Range("H89").Value = 14476.75
$ws.Range("I89").Value = 2133.3333
$ws.Range("J89").Value = 21882.8
$ws.Range("K89").Value = 10666.6665
$ws.Range("L89").Value = 109414
$ws.Range("M89").Value = -5050.666499999999
$ws.Range("N89").Value = -120646

$ws.Range("H105").Value = 11113718
$ws.Range("I105").Value = 13335881
$ws.Range("J105").Value = 2902.2
$ws.Range("K105").Value = 13335881
$ws.Range("L105").Value = 2902.2
$ws.Range("M105").Value = -13334134
$ws.Range("N105").Value = -6396.2

$ws.Range("H134").Value = 5452.4165
$ws.Range("I134").Value = 4959.6
$ws.Range("K134").Value = 14878.8
$ws.Range("M134").Value = -12343.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()

$ws.Range("H134").Value = 3062.0454
$ws.Range("I134").Value = 1415.6666
$ws.Range("J134").Value = 6590
$ws.Range("K134").Value = 4246.9998
$ws.Range("L134").Value = 19770
$ws.Range("M134").Value = -1711.9998
$ws.Range("N134").Value = -24840

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1039.23
$ws.Range("I68").Value = 786.05084
$ws.Range("J68").Value = 1403.561
$ws.Range("K68").Value = 2358.15252
$ws.Range("L68").Value = 4210.683
$ws.Range("M68").Value = -1547.15252
$ws.Range("N68").Value = -5832.683

$ws.Range("H71").Value = 1039.23
$ws.Range("I71").Value = 786.05084
$ws.Range("J71").Value = 1403.561
$ws.Range("K71").Value = 7074.45756
$ws.Range("L71").Value = 12632.049
$ws.Range("M71").Value = -3018.45756
$ws.Range("N71").Value = -20744.049

$ws.Range("H131").Value = 2719.863
$ws.Range("I131").Value = 925
$ws.Range("J131").Value = 2823.913
$ws.Range("K131").Value = 2775
$ws.Range("L131").Value = 8471.739
$ws.Range("M131").Value = 2265
$ws.Range("N131").Value = -18551.739

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6463.1577
$ws.Range("I70").Value = 6862.5
$ws.Range("J70").Value = 4333.3335
$ws.Range("K70").Value = 6862.5
$ws.Range("L70").Value = 4333.3335
$ws.Range("M70").Value = -6592.5
$ws.Range("N70").Value = -4873.3335

$ws.Range("H73").Value = 6463.1577
$ws.Range("I73").Value = 6862.5
$ws.Range("J73").Value = 4333.3335
$ws.Range("K73").Value = 6862.5
$ws.Range("L73").Value = 4333.3335
$ws.Range("M73").Value = -5926.5
$ws.Range("N73").Value = -6205.3335

$ws.Range("H126").Value = 2447.3157
$ws.Range("I126").Value = 2333.1667
$ws.Range("K126").Value = 6999.500100000001
$ws.Range("M126").Value = -4529.500100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 32468
$ws.Range("J25").Value = 14936
$ws.Range("L25").Value = 14936
$ws.Range("N25").Value = -15396

$ws.Range("H40").Value = 4503.2354
$ws.Range("I40").Value = 2740.6667
$ws.Range("J40").Value = 4880.9287
$ws.Range("K40").Value = 2740.6667
$ws.Range("L40").Value = 4880.9287
$ws.Range("M40").Value = -2604.6667
$ws.Range("N40").Value = -5152.9287

$ws.Range("H69").Value = 33000
$ws.Range("J69").Value = 33000
$ws.Range("L69").Value = 33000
$ws.Range("N69").Value = -34622

$ws.Range("H72").Value = 33000
$ws.Range("J72").Value = 33000
$ws.Range("L72").Value = 99000
$ws.Range("N72").Value = -107112

$ws.Range("H122").Value = 3990.9092
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 3990.9092
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 11972.7276
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -16872.7276

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 14287693
$ws.Range("I132").Value = 21740716
$ws.Range("K132").Value = 65222148
$ws.Range("M132").Value = -65219618
